# Commit: "Suppression des references et correctifs 96ac5840dfc0f40661f77a6732bdab5ef2bee7f9"
#
# Applies the changes:
#  1. Metadata sheet (B8): update the "Date" property value to the new
#     regeneration timestamp.
#  2. Elements sheet, row 7 (ActorSNR.XCN9.composant1): Min/Max (F7/G7) and
#     Base Min/Base Max (AG7/AH7) change from "1" to "0". A leading
#     apostrophe is used so these numeric-looking values are stored as text
#     (matching the source data, which keeps Min/Max as text labels) rather
#     than being auto-converted to numbers.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-05-05T11:54:16+00:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("F7").Value = "'0"
$elements.Range("G7").Value = "'0"
$elements.Range("AG7").Value = "'0"
$elements.Range("AH7").Value = "'0"
